# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4981.294
$ws.Range("I28").Value = 6235.5
$ws.Range("K28").Value = 6235.5
$ws.Range("M28").Value = -5750.5
$ws.Range("H40").Value = 3591.85
$ws.Range("I40").Value = 3129.2
$ws.Range("K40").Value = 3129.2
$ws.Range("M40").Value = -2954.2
$ws.Range("H62").Value = 4763.143
$ws.Range("I62").Value = 4744.923
$ws.Range("K62").Value = 4744.923
$ws.Range("M62").Value = -4120.923
$ws.Range("H65").Value = 4763.143
$ws.Range("I65").Value = 4744.923
$ws.Range("K65").Value = 23724.615
$ws.Range("M65").Value = -20604.615
$ws.Range("H112").Value = 2723.9167
$ws.Range("I112").Value = 1341.3334
$ws.Range("J112").Value = 3184.7778
$ws.Range("K112").Value = 4024.0002
$ws.Range("L112").Value = 9554.3334
$ws.Range("M112").Value = -2916.0002
$ws.Range("N112").Value = -11770.3334
$ws.Range("H116").Value = 11114948
$ws.Range("I116").Value = 15280942
$ws.Range("J116").Value = 5631.6665
$ws.Range("K116").Value = 15280942
$ws.Range("L116").Value = 5631.6665
$ws.Range("M116").Value = -15277500
$ws.Range("N116").Value = -12515.6665
$ws.Range("H134").Value = 89993.39999999999
$ws.Range("J134").Value = 89993.39999999999
$ws.Range("L134").Value = 89993.39999999999
$ws.Range("N134").Value = -100133.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4177.6807
$ws.Range("I32").Value = 3918.9333
$ws.Range("K32").Value = 3918.9333
$ws.Range("M32").Value = -3631.9333
$ws.Range("H61").Value = 5998.6665
$ws.Range("I61").Value = 4998
$ws.Range("K61").Value = 4998
$ws.Range("M61").Value = -4786
$ws.Range("H102").Value = 13143.04
$ws.Range("I102").Value = 27598.75
$ws.Range("J102").Value = 6340.353
$ws.Range("K102").Value = 27598.75
$ws.Range("L102").Value = 6340.353
$ws.Range("M102").Value = -25976.75
$ws.Range("N102").Value = -9584.352999999999
$ws.Range("H122").Value = 762043.3
$ws.Range("I122").Value = 2802.805
$ws.Range("K122").Value = 8408.414999999999
$ws.Range("M122").Value = -5958.414999999999
$ws.Range("H136").Value = 5998.6665
$ws.Range("I136").Value = 4998
$ws.Range("K136").Value = 14994
$ws.Range("M136").Value = -12444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 11044
$ws.Range("I94").Value = 13637.904
$ws.Range("J94").Value = 1965.3334
$ws.Range("K94").Value = 13637.904
$ws.Range("L94").Value = 1965.3334
$ws.Range("M94").Value = -13186.904
$ws.Range("N94").Value = -2867.3334
$ws.Range("H99").Value = 13565.637
$ws.Range("I99").Value = 15736.12
$ws.Range("J99").Value = 6782.875
$ws.Range("K99").Value = 15736.12
$ws.Range("L99").Value = 6782.875
$ws.Range("M99").Value = -14238.12
$ws.Range("N99").Value = -9778.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6647.971
$ws.Range("I31").Value = 7921.8096
$ws.Range("J31").Value = 4737.2144
$ws.Range("K31").Value = 7921.8096
$ws.Range("L31").Value = 4737.2144
$ws.Range("M31").Value = -7626.8096
$ws.Range("N31").Value = -5327.2144
$ws.Range("H34").Value = 6647.971
$ws.Range("I34").Value = 7921.8096
$ws.Range("J34").Value = 4737.2144
$ws.Range("K34").Value = 7921.8096
$ws.Range("L34").Value = 4737.2144
$ws.Range("M34").Value = -7719.8096
$ws.Range("N34").Value = -5141.2144
$ws.Range("H62").Value = 64824.547
$ws.Range("J62").Value = 95317.14
$ws.Range("L62").Value = 95317.14
$ws.Range("N62").Value = -96565.14
$ws.Range("H65").Value = 64824.547
$ws.Range("J65").Value = 95317.14
$ws.Range("L65").Value = 476585.7
$ws.Range("N65").Value = -482825.7
$ws.Range("H86").Value = 9086.916999999999
$ws.Range("I86").Value = 7341
$ws.Range("J86").Value = 10832.833
$ws.Range("K86").Value = 7341
$ws.Range("L86").Value = 10832.833
$ws.Range("M86").Value = -6218
$ws.Range("N86").Value = -13078.833
$ws.Range("H89").Value = 9086.916999999999
$ws.Range("I89").Value = 7341
$ws.Range("J89").Value = 10832.833
$ws.Range("K89").Value = 36705
$ws.Range("L89").Value = 54164.165
$ws.Range("M89").Value = -31089
$ws.Range("N89").Value = -65396.165
$ws.Range("H132").Value = 1821.4736
$ws.Range("I132").Value = 1706
$ws.Range("K132").Value = 5118
$ws.Range("M132").Value = -2588
$ws.Range("H134").Value = 5405.4287
$ws.Range("I134").Value = 7460
$ws.Range("J134").Value = 2666
$ws.Range("K134").Value = 22380
$ws.Range("L134").Value = 7998
$ws.Range("M134").Value = -19845
$ws.Range("N134").Value = -13068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 125.066666
$ws.Range("I12").Value = 316.4
$ws.Range("K12").Value = 949.1999999999999
$ws.Range("M12").Value = -776.1999999999999
$ws.Range("H103").Value = 983.73334
$ws.Range("J103").Value = 1512.2222
$ws.Range("L103").Value = 4536.6666
$ws.Range("N103").Value = -6294.6666
$ws.Range("H115").Value = 3499
$ws.Range("I115").Value = 550
$ws.Range("J115").Value = 4973.5
$ws.Range("K115").Value = 1650
$ws.Range("L115").Value = 14920.5
$ws.Range("M115").Value = -475
$ws.Range("N115").Value = -17270.5
$ws.Range("H131").Value = 1610.3489
$ws.Range("J131").Value = 1633.6962
$ws.Range("L131").Value = 4901.0886
$ws.Range("N131").Value = -14981.0886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 172.90909
$ws.Range("I2").Value = 155.77777
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 155.77777
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -42.77777
$ws.Range("N2").Value = -476
$ws.Range("H70").Value = 9607.571
$ws.Range("I70").Value = 9060.125
$ws.Range("J70").Value = 10337.5
$ws.Range("K70").Value = 9060.125
$ws.Range("L70").Value = 10337.5
$ws.Range("M70").Value = -8790.125
$ws.Range("N70").Value = -10877.5
$ws.Range("H73").Value = 9607.571
$ws.Range("I73").Value = 9060.125
$ws.Range("J73").Value = 10337.5
$ws.Range("K73").Value = 9060.125
$ws.Range("L73").Value = 10337.5
$ws.Range("M73").Value = -8124.125
$ws.Range("N73").Value = -12209.5
$ws.Range("H113").Value = 14970.333
$ws.Range("I113").Value = 23606.6
$ws.Range("K113").Value = 23606.6
$ws.Range("M113").Value = -21436.6
$ws.Range("H122").Value = 6735.4736
$ws.Range("I122").Value = 4111.1934
$ws.Range("K122").Value = 12333.5802
$ws.Range("M122").Value = -9883.5802
$ws.Range("H132").Value = 1999.6666
$ws.Range("I132").Value = 1999.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5998.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3468.5
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 7553
$ws.Range("J43").Value = 7553
$ws.Range("L43").Value = 7553
$ws.Range("N43").Value = -7939
$ws.Range("H93").Value = 8725.23
$ws.Range("I93").Value = 9192.9
$ws.Range("J93").Value = 7166.3335
$ws.Range("K93").Value = 9192.9
$ws.Range("L93").Value = 7166.3335
$ws.Range("M93").Value = -7944.9
$ws.Range("N93").Value = -9662.333500000001
$ws.Range("H122").Value = 6824.6875
$ws.Range("I122").Value = 6799.4443
$ws.Range("K122").Value = 20398.3329
$ws.Range("M122").Value = -17948.3329
$ws.Range("H132").Value = 555170.25
$ws.Range("I132").Value = 878182.75
$ws.Range("J132").Value = 6048.9
$ws.Range("K132").Value = 2634548.25
$ws.Range("L132").Value = 18146.7
$ws.Range("M132").Value = -2632018.25
$ws.Range("N132").Value = -23206.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 11866.25
$ws.Range("I107").Value = 2833.2964
$ws.Range("J107").Value = 60644.2
$ws.Range("K107").Value = 8499.889200000001
$ws.Range("L107").Value = 181932.6
$ws.Range("M107").Value = -6579.889200000001
$ws.Range("N107").Value = -185772.6
$ws.Range("H122").Value = 5357.7827
$ws.Range("I122").Value = 2384.5454
$ws.Range("J122").Value = 8083.25
$ws.Range("K122").Value = 7153.6362
$ws.Range("L122").Value = 24249.75
$ws.Range("M122").Value = -4703.6362
$ws.Range("N122").Value = -29149.75
$ws.Range("H126").Value = 29884.867
$ws.Range("I126").Value = 51419.75
$ws.Range("J126").Value = 5273.5713
$ws.Range("K126").Value = 154259.25
$ws.Range("L126").Value = 15820.7139
$ws.Range("M126").Value = -151789.25
$ws.Range("N126").Value = -20760.7139
$ws.Range("H132").Value = 15291.272
$ws.Range("I132").Value = 17244.973
$ws.Range("K132").Value = 51734.91900000001
$ws.Range("M132").Value = -49204.91900000001
